$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.095157742500305
$ws.Range("B1").Value = 2.995087623596191
$ws.Range("C1").Value = 6.420498847961426
$ws.Range("D1").Value = 4.179830074310303
$ws.Range("E1").Value = 1.345168828964233
